$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.92"
$ws.Range("E2").Value = "'-0.81%"
$ws.Range("D3").Value = "'36.89"
$ws.Range("E3").Value = "'3.55%"
$ws.Range("D4").Value = "'5.010"
$ws.Range("E4").Value = "'-0.91%"
$ws.Range("D5").Value = "'0.07675"
$ws.Range("E5").Value = "'-1.61%"
$ws.Range("D6").Value = "'2.058"
$ws.Range("E6").Value = "'-8.69%"
$ws.Range("D7").Value = "'7.969"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.033"
$ws.Range("E8").Value = "'-0.34%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9158"
$ws.Range("E9").Value = "'-1.50%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09462"
$ws.Range("E10").Value = "'-0.27%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1852"
$ws.Range("E11").Value = "'1.44%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08591"
$ws.Range("E12").Value = "'0.35%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03602"
$ws.Range("E13").Value = "'5.20%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09998"
$ws.Range("E14").Value = "'0.53%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001475"
$ws.Range("E15").Value = "'-0.32%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005738"
$ws.Range("E16").Value = "'-0.18%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.473"
$ws.Range("E17").Value = "'-0.22%"
$ws.Range("D18").Value = "'2.380"
$ws.Range("E18").Value = "'9.21%"
$ws.Range("D19").Value = "'0.3368"
$ws.Range("E19").Value = "'-1.12%"
$ws.Range("D20").Value = "'0.1332"
$ws.Range("E20").Value = "'0.67%"
$ws.Range("D21").Value = "'4.956"
$ws.Range("E21").Value = "'9.07%"
$ws.Range("D22").Value = "'0.2213"
$ws.Range("E22").Value = "'-1.07%"
$ws.Range("D23").Value = "'0.04609"
$ws.Range("E23").Value = "'-1.62%"
$ws.Range("E24").Value = "'11.98%"
$ws.Range("D25").Value = "'0.001241"
$ws.Range("E25").Value = "'0.00%"
$ws.Range("D26").Value = "'0.0001412"
$ws.Range("E26").Value = "'8.54%"
$ws.Range("D39").Value = "'0.01745"
$ws.Range("E39").Value = "'-1.73%"
$ws.Range("D40").Value = "'0.04591"
$ws.Range("E40").Value = "'-2.68%"
$ws.Range("D41").Value = "'0.007712"
$ws.Range("E41").Value = "'-2.21%"
$ws.Range("D42").Value = "'0.1390"
$ws.Range("E42").Value = "'-2.03%"
$ws.Range("D43").Value = "'0.007783"
$ws.Range("E43").Value = "'-2.82%"
$ws.Range("D44").Value = "'0.002172"
$ws.Range("E44").Value = "'-2.30%"
$ws.Range("D45").Value = "'0.01036"
$ws.Range("E45").Value = "'14.02%"
$ws.Range("D46").Value = "'0.00006308"
$ws.Range("E46").Value = "'1.87%"
$ws.Range("E47").Value = "'0.72%"
$ws.Range("D48").Value = "'0.0005834"
$ws.Range("E48").Value = "'0.58%"
$ws.Range("D49").Value = "'34.63"
$ws.Range("E49").Value = "'557.61%"
$ws.Range("E50").Value = "'-25.11%"
$ws.Range("E51").Value = "'0.72%"
